$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.178.94"
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("D3").Value = "'2.471.27"
$ws.Range("E3").Value = "  +0.06%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'583.97"
$ws.Range("E5").Value = "  +0.19%  "
$ws.Range("D6").Value = "'174.57"
$ws.Range("E6").Value = "  +3.74%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'0.513"
$ws.Range("E8").Value = "  -0.19%  "
$ws.Range("E9").Value = "  +2.60%  "
$ws.Range("E10").Value = "  +0.53%  "
$ws.Range("D11").Value = "'4.95"
$ws.Range("E11").Value = "  -0.38%  "
$ws.Range("D12").Value = "'0.332"
$ws.Range("E12").Value = "  +0.59%  "
$ws.Range("D13").Value = "'2.932.79"
$ws.Range("E13").Value = "  +0.59%  "
$ws.Range("D14").Value = "'25.45"
$ws.Range("E14").Value = "  -0.37%  "
$ws.Range("D15").Value = "'67.073.57"
$ws.Range("E15").Value = "  +0.04%  "
$ws.Range("E16").Value = "  +0.08%  "
$ws.Range("D17").Value = "'2.578.55"
$ws.Range("E17").Value = "  +2.01%  "
$ws.Range("D18").Value = "'7.52"
$ws.Range("E18").Value = "  -1.33%  "
$ws.Range("D19").Value = "'10.93"
$ws.Range("E19").Value = "  -2.18%  "
$ws.Range("D20").Value = "'350.37"
$ws.Range("E20").Value = "  -0.94%  "
$ws.Range("D21").Value = "'3.99"
$ws.Range("E21").Value = "  -1.26%  "
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").Value = "'69.22"
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("D24").Value = "'4.22"
$ws.Range("E24").Value = "  -0.28%  "
$ws.Range("E25").Value = "  +1.16%  "
$ws.Range("D26").Value = "'9.16"
$ws.Range("E26").Value = "  -0.93%  "
$ws.Range("E27").Value = "  +0.45%  "
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("D29").Value = "'0.0₃0903"
$ws.Range("E29").Value = "  -0.35%  "
$ws.Range("D30").Value = "'500.81"
$ws.Range("E30").Value = "  -2.58%  "
$ws.Range("D31").Value = "'7.75"
$ws.Range("E31").Value = "  +0.27%  "
$ws.Range("D32").Value = "'1.23"
$ws.Range("E32").Value = "  -0.51%  "
$ws.Range("E33").Value = "  -1.18%  "
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("D36").Value = "'162.04"
$ws.Range("E36").Value = "  +1.64%  "
$ws.Range("D37").Value = "'18.68"
$ws.Range("E37").Value = "  +0.08%  "
$ws.Range("D38").Value = "'18.13"
$ws.Range("E38").Value = "  -1.35%  "
$ws.Range("E39").Value = "  -1.36%  "
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("E41").Value = "  +1.81%  "
$ws.Range("D42").Value = "'0.328"
$ws.Range("E42").Value = "  +0.26%  "
$ws.Range("D43").Value = "'4.83"
$ws.Range("E43").Value = "  +0.81%  "
$ws.Range("D44").Value = "'2.39"
$ws.Range("E44").Value = "  +1.87%  "
$ws.Range("D45").Value = "'142.04"
$ws.Range("E45").Value = "  +0.83%  "
$ws.Range("D46").Value = "'3.48"
$ws.Range("E46").Value = "  +0.58%  "
$ws.Range("D47").Value = "'0.513"
$ws.Range("E47").Value = "  -0.47%  "
$ws.Range("D48").Value = "'0.0₆0255"
$ws.Range("E48").Value = "  +1.51%  "
$ws.Range("D49").Value = "'0.0738"
$ws.Range("E49").Value = "  +0.26%  "
$ws.Range("E50").Value = "  -0.87%  "
$ws.Range("D51").Value = "'0.582"
$ws.Range("E51").Value = "  +0.31%  "
